$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.924774646759033
$ws.Range("B1").Value = 5.862129211425781
$ws.Range("C1").Value = 6.59728479385376
$ws.Range("D1").Value = 7.048226356506348
$ws.Range("E1").Value = 4.171567916870117
